$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 170, shifting the existing rows 170-190 down to 171-191.
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new weekly data point.
$ws.Cells.Item(170, 1).Value = 8
$ws.Cells.Item(170, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = 45212
$ws.Cells.Item(170, 5).Value = 4
$ws.Cells.Item(170, 6).Value = 100114007
$ws.Cells.Item(170, 7).Value = "Jengibre"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 360
$ws.Cells.Item(170, 11).Value = 21000
$ws.Cells.Item(170, 12).Value = 22000
$ws.Cells.Item(170, 13).Value = 21500
$ws.Cells.Item(170, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(170, 15).Value = "Perú"
$ws.Cells.Item(170, 16).Value = 1654
$ws.Cells.Item(170, 17).Value = 13
$ws.Cells.Item(170, 18).Value = "Hortaliza"
